$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert to the sprint-II checklist state: drop the last 3 appended rows
# (rows 24-26 -- the "Tshering Dorji" / "Tshering Norbu" sightings), shrinking
# the data range back from A1:K26 to A1:K23.
$ws.Rows.Item(24).Resize(3).Delete()

# Restore column E's width (used by the sprint-II analysis pass) and the
# previous cell selection.
$ws.Columns.Item(5).ColumnWidth = 10.6
$ws.Range("E6").Select()
